$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data set gained a new (most recent) weekly record. It is inserted at
# row 5 (the data rows start at row 2; row 5 is the first "Provincia de
# Quillota" origin record), pushing the existing rows 5-51 down to 6-52.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly record.
$ws.Cells.Item(5, 1).Value = 9
$ws.Cells.Item(5, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(5, 3).Value = "Metropolitana"
$ws.Cells.Item(5, 4).Value = 45022
$ws.Cells.Item(5, 5).Value = 13
$ws.Cells.Item(5, 6).Value = 100112010
$ws.Cells.Item(5, 7).Value = "Achicoria"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 60
$ws.Cells.Item(5, 11).Value = 7000
$ws.Cells.Item(5, 12).Value = 7000
$ws.Cells.Item(5, 13).Value = 7000
$ws.Cells.Item(5, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(5, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(5, 16).Value = 438
$ws.Cells.Item(5, 17).Value = 16
$ws.Cells.Item(5, 18).Value = "Hortaliza"
